# Rename the "Address" column-name property to "Bill_address" in the
# last table of the document (the row whose Description cell reads
# "Bill address"), per the commit: "Change attribute name of property".
$d = $word.ActiveDocument

$table = $d.Tables.Item(9)
$cell = $table.Cell(4, 1)
$para = $cell.Range.Paragraphs.Item(1)
$paraRange = $para.Range

# The cell currently holds two runs: "A" (with an eastAsia rFonts hint)
# followed by "ddress" (no rFonts hint). Delete the leading "A" run's
# single character, then prepend "Bill_a" onto the remaining "ddress"
# text so the final text becomes "Bill_address" as one run that keeps
# the plain (non-hinted) formatting of the original "ddress" run.
$start = $paraRange.Start
$firstChar = $d.Range($start, $start + 1)
$firstChar.Delete()

$remaining = $para.Range
$target = $d.Range($remaining.Start, $remaining.Start + 1)
$target.InsertBefore("Bill_a")

$verifyCell = $d.Tables.Item(9).Cell(4, 1)
Write-Output "CellText=[$($verifyCell.Range.Text)]"
